$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header cell A1: "No." (already "No." textually, just re-assert)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "No."

# ---------------------------------------------------------------------
# 2. Row 2 - "Test add to-do item without entering a value"
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Test add to-do item without entering a value"
$ws.Range("C2").Value = "To check if to-do page has error validation when the user attempts to add an empty to-do item."
$ws.Range("E2").Value = "There will be an error prompt to enter a value in the to-do item field."
$ws.Range("F2").Value = "Same as expected outcome."
$ws.Range("G2").Value = "Pass"

# ---------------------------------------------------------------------
# 3. Row 3 - "Test adding a to-do item"
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "To check if a to-do item can be added successfully in the to-do page."
$ws.Range("E3").Value = "The to-do item will be added successfully and be displayed in the to-do page."
$ws.Range("F3").Value = "Same as expected outcome."
$ws.Range("G3").Value = "Pass"

# ---------------------------------------------------------------------
# 4. Row 4 - "Test deleting a to-do item"
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "To check if a to-do item can be deleted successfully in the to-do page."
$ws.Range("E4").Value = "The to-do item will be deleted successfully and will not be shown in the to-do page."
$ws.Range("F4").Value = "Same as expected outcome."
$ws.Range("G4").Value = "Pass"

# ---------------------------------------------------------------------
# 5. Row 5 - now "Test automatically archiving after completing a to-do item"
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "Test automatically archiving after completing a to-do item"
$ws.Range("C5").Value = "To check if a to-do item can be archived successfully in the to-do page."
$ws.Range("E5").Value = "Item is archived successfully and will not be shown in the to-do page."
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "Fail"

# ---------------------------------------------------------------------
# 6. Row 6 - "Test navigation to To-do Page"
# ---------------------------------------------------------------------
$ws.Range("C6").Value = "To test if the to-do page can be navigated to successfully."
$ws.Range("E6").Value = "The user will be brought to the to-do page."
$ws.Range("F6").Value = "Same as expected outcome."
$ws.Range("G6").Value = "Pass"

# ---------------------------------------------------------------------
# 7. New rows 7-9: copy formatting from row 6 (style index 5, bordered)
#    then fill in the new test case data.
# ---------------------------------------------------------------------
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A7:G9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Row 7 - Test successful display of to-do items based on user
$ws.Range("A7").Value = 2.5
$ws.Range("B7").Value = "Test successful display of to-do items based on user"
$ws.Range("C7").Value = "To check if to-do items can only be seen by the user that has created it."
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "The to-do items that are unique to the user will be shown successfully."
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = "Fail"
$ws.Rows.Item(7).RowHeight = 51.75

# Row 8 - Test displaying of timestamp (creation date) of to-do item
$ws.Range("A8").Value = 2.6
$ws.Range("B8").Value = "Test displaying of timestamp (creation date) of to-do item"
$ws.Range("C8").Value = "To check if to-do items have individual timestamps of creation datetime."
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "The to-do items with timestamps will be shown successfully."
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "Fail"
$ws.Rows.Item(8).RowHeight = 72.75

# Row 9 - Test displaying of individual to-do items for each user
$ws.Range("A9").Value = 2.7
$ws.Range("B9").Value = "Test displaying of individual to-do items for each user"
$ws.Range("C9").Value = "To test if the to-do items shown are the ones created by the user and not other users."
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "After user log in, he/she sees only his/her own to-do items. Logging into another user's account will show a different list of to-do items."
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "Fail"
$ws.Rows.Item(9).RowHeight = 89.25

# ---------------------------------------------------------------------
# 8. Row 10 used to hold (now-relocated) data; it becomes a blank,
#    unstyled row again.
# ---------------------------------------------------------------------
$ws.Range("A10:G10").Clear() | Out-Null
$ws.Rows.Item(10).RowHeight = 50.25

# ---------------------------------------------------------------------
# 9. Sheet view: zoom to 85%, and move the active selection to B5.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B5").Select() | Out-Null
